# Add a "Repo.zipname" column: the specimen photos/spectra for each genus
# are zipped together for storage in the repo, named "<Family>_<Genus>.zip".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$familyCol = 2   # column B
$genusCol = 3    # column C
$zipCol = 21     # column U (one past the last existing column, T)

$ws.Cells.Item(1, $zipCol).Value = "Repo.zipname"

for ($r = 2; $r -le $lastRow; $r++) {
    $family = $ws.Cells.Item($r, $familyCol).Text
    $genus = $ws.Cells.Item($r, $genusCol).Text
    $ws.Cells.Item($r, $zipCol).Value = $family + "_" + $genus + ".zip"
}
